$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to be treated as text so numeric-looking strings
# (e.g. "1.000", "0.9991", "0.00000000118") are preserved exactly
# instead of being silently re-interpreted as numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.484.41'
$ws.Range("E2").Value = '  -0.98%  '
$ws.Range("D3").Value = '1.849.27'
$ws.Range("E3").Value = '  -0.55%  '
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '241.75'
$ws.Range("E5").Value = '  -1.12%  '
$ws.Range("D6").Value = '0.6289'
$ws.Range("E6").Value = '  -2.28%  '
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '48.05'
$ws.Range("E8").Value = '  -0.15%  '
$ws.Range("D9").Value = '0.07535'
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("D10").Value = '0.2976'
$ws.Range("E10").Value = '  -0.31%  '
$ws.Range("D11").Value = '24.34'
$ws.Range("E11").Value = '  -1.16%  '
$ws.Range("D12").Value = '0.07702'
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("D13").Value = '1.917.15'
$ws.Range("E13").Value = '  +3.07%  '
$ws.Range("E14").Value = '  -0.86%  '
$ws.Range("D15").Value = '0.6883'
$ws.Range("E15").Value = '  -0.54%  '
$ws.Range("D16").Value = '83.71'
$ws.Range("E16").Value = '  -0.38%  '
$ws.Range("D17").Value = '0.000009811'
$ws.Range("E17").Value = '  -0.99%  '
$ws.Range("D18").Value = '2.163.07'
$ws.Range("E18").Value = '  +2.36%  '
$ws.Range("D19").Value = '6.249'
$ws.Range("E19").Value = '  +1.84%  '
$ws.Range("D20").Value = '29.566.13'
$ws.Range("E20").Value = '  -0.78%  '
$ws.Range("D21").Value = '233.77'
$ws.Range("D22").Value = '12.49'
$ws.Range("E22").Value = '  -1.33%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").Value = '7.633'
$ws.Range("E24").Value = '  +0.91%  '
$ws.Range("D25").Value = '1.000'
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").Value = '154.49'
$ws.Range("E26").Value = '  -2.58%  '
$ws.Range("D27").Value = '0.1392'
$ws.Range("E27").Value = '  -2.18%  '
$ws.Range("D28").Value = '8.450'
$ws.Range("E28").Value = '  -1.39%  '
$ws.Range("D29").Value = '17.72'
$ws.Range("E29").Value = '  -1.19%  '
$ws.Range("D30").Value = '1.477'
$ws.Range("E30").Value = '  -1.20%  '
$ws.Range("D31").Value = '0.05849'
$ws.Range("E31").Value = '  -6.16%  '
$ws.Range("D32").Value = '1.256'
$ws.Range("E32").Value = '  -2.65%  '
$ws.Range("D33").Value = '4.104'
$ws.Range("E33").Value = '  -1.27%  '
$ws.Range("D34").Value = '4.039'
$ws.Range("E34").Value = '  -1.44%  '
$ws.Range("D35").Value = '1.882'
$ws.Range("E35").Value = '  -0.97%  '
$ws.Range("E36").Value = '  -0.39%  '
$ws.Range("D37").Value = '0.7193'
$ws.Range("E37").Value = '  -1.79%  '
$ws.Range("D38").Value = '2.587'
$ws.Range("E38").Value = '  -0.85%  '
$ws.Range("D39").Value = '1.243.47'
$ws.Range("E39").Value = '  +1.69%  '
$ws.Range("D40").Value = '2.799'
$ws.Range("E40").Value = '  -0.83%  '
$ws.Range("D41").Value = '0.01784'
$ws.Range("E41").Value = '  -0.33%  '
$ws.Range("D42").Value = '0.9058'
$ws.Range("E42").Value = '  -1.55%  '
$ws.Range("D43").Value = '6.159'
$ws.Range("E43").Value = '  -2.41%  '
$ws.Range("D44").Value = '2.068.40'
$ws.Range("E44").Value = '  +1.99%  '
$ws.Range("D45").Value = '0.9997'
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("D46").Value = '102.03'
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").Value = '67.25'
$ws.Range("E47").Value = '  +0.18%  '
$ws.Range("D48").Value = '7.312'
$ws.Range("E48").Value = '  +8.48%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.00000000118'
$ws.Range("E49").Value = '  -0.43%  '
$ws.Range("D50").Value = '0.4042'
$ws.Range("E50").Value = '  -0.66%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '9.159'
$ws.Range("E51").Value = '  -0.41%  '
